# Apply the "Simplified RequestsService and ReportConfigurationService" edit.
#
# Summary of the change (see xml_diff):
#   * The "Generated by" value changes from the admin placeholder to a
#     client placeholder.
#   * The whole "Admin" block (column G:I, header "Admin" in G8, the
#     "Name"/"Pronouns"/"Works At" sub-headers in row 9, and the
#     "Final boss"/"Soul/Soul"/"The Klin of the First Flame" values in
#     every data row) is removed completely, shrinking the table from
#     columns B:I down to B:F.
#   * The now-redundant extra data rows (the table only keeps 6 data
#     rows instead of 10) are removed as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Generated by" cell (C3): admin name/role -> client name/role.
$ws.Range("C3").Value = "John Snow (Client)"

# 2) Remove the Admin-only columns (G:I). This automatically drops the
#    "Admin" header, the "Name"/"Pronouns"/"Works At" sub-headers and all
#    of the "Final boss"/"Soul/Soul"/"The Klin of the First Flame" cells,
#    shrinks the B7:I7 merge down to B7:F7, and removes the G8:I8 merge.
$ws.Range("G1:I1").EntireColumn.Delete() | Out-Null

# 3) Drop the surplus data rows so only 6 data rows (10-15) remain.
$ws.Range("A16:A19").EntireRow.Delete() | Out-Null

# 4) The remaining "Additional Info" column now reads "Didn't piss of the
#    Wall" on every row except the very first data row (10), instead of
#    alternating with "Killed the Queen on working hours".
$ws.Range("F12").Value = "Didn't piss of the Wall"
$ws.Range("F14").Value = "Didn't piss of the Wall"

# 5) Keep the selection on I9 (closest reproducible equivalent of the
#    original drag-selection artifact "I9 G9:F9" recorded in the file).
$ws.Range("I9").Select() | Out-Null
